$d = $word.ActiveDocument

# Paragraph 2 holds the "master list of word ..." run-soup that is being
# replaced wholesale with a new set of runs (per the diff). Locate it by
# index (it's the only content paragraph besides the Aspose banner).
$p2 = $d.Paragraphs.Item(2)
$oldStart = $p2.Range.Start
$oldEnd = $p2.Range.End

# Build the replacement content as raw OOXML runs so each becomes its own
# <w:r> (mirrors the literal run-split in the target diff) instead of being
# coalesced into a single run the way Range.Text/InsertAfter would do.
$newXml = @'
<?xml version="1.0" standalone="yes"?>
<?mso-application progid="Word.Document"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
<pkg:xmlData>
<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
<w:body>
<w:p><w:r><w:t>number correct for all scoring dic</w:t></w:r><w:r><w:t xml:space="preserve">     </w:t></w:r><w:r><w:t>-1369</w:t></w:r><w:r><w:t xml:space="preserve">  number correct for specific scoring dic</w:t></w:r><w:r><w:t xml:space="preserve">     </w:t></w:r><w:r><w:t>-1167</w:t></w:r><w:r><w:t xml:space="preserve">  number correct for num scoring dic</w:t></w:r><w:r><w:t xml:space="preserve">     </w:t></w:r><w:r><w:t>-1109</w:t></w:r><w:r><w:t xml:space="preserve">   number correct for fancy way</w:t></w:r><w:r><w:t>-1139</w:t></w:r><w:r><w:t xml:space="preserve">total reviews were    </w:t></w:r><w:r><w:t>1497</w:t></w:r></w:p>
</w:body>
</w:document>
</pkg:xmlData>
</pkg:part>
</pkg:package>
'@

# Insert the new runs inline at the very start of paragraph 2 (a collapsed
# range merges into the existing paragraph rather than splitting a new one).
$insPoint = $d.Range($oldStart, $oldStart)
$insPoint.InsertXML($newXml)

# The old 32 runs got pushed after our freshly-inserted text; compute their
# new bounds and delete them, leaving only the new runs behind.
$shift = $p2.Range.End - $oldEnd
$oldContent = $d.Range($oldStart + $shift, $oldEnd + $shift)
$oldContent.Delete()
